$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header cells: _old -> _FV2410, _new -> _FV2504 ---
$headers = @(
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
  "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410", "Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504",
  "Segment ID_FV2504", "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Freeze the header row (split below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn the data range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), $null, 1)
$tbl.Name = "Table1"
